$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (ano 2025) metrics: total_customers and new_customers changed,
# which in turn changes the derived new_rate / returning_rate percentages.
$ws.Range("C6").Value = 396
$ws.Range("E6").Value = 91
$ws.Range("G6").Value = 22.97979797979798
$ws.Range("H6").Value = 77.02020202020202
